$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "A57" = "Aprosmictus erythropterus"
    "B57" = "https://www.xeno-canto.org/354254/download"
    "C57" = "Bowra Station, Queensland, Australia"
    "D57" = "Greg McLachlan"
    "E57" = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
    "F57" = "assets/misc/cc.png"
    "G57" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A58" = "Apus pacificus"
    "B58" = "https://www.xeno-canto.org/352598/download"
    "C58" = "Cattana Wetlands, Queensland, Australia"
    "D58" = "Marc Anderson"
    "E58" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F58" = "assets/misc/cc.png"
    "G58" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A59" = "Apus pacificus"
    "B59" = "https://www.xeno-canto.org/286311/download"
    "C59" = "Kiritappu, Hokkaido, Japan"
    "D59" = "Peter Boesman"
    "E59" = "https://www.xeno-canto.org/contributor/OOECIWCSWV"
    "F59" = "assets/misc/cc.png"
    "G59" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A60" = "Ardea alba"
    "B60" = "https://www.xeno-canto.org/365950/download"
    "C60" = "Freemans Reach, New South Wales, Australia"
    "D60" = "Greg McLachlan"
    "E60" = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
    "F60" = "assets/misc/cc.png"
    "G60" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A61" = "Ardea alba"
    "B61" = "https://www.xeno-canto.org/431288/download"
    "C61" = "Reserva Natural Palmarí, Rio Javarí, Brazil"
    "D61" = "Jerome Fischer"
    "E61" = "https://www.xeno-canto.org/contributor/JPBSNBUUEF"
    "F61" = "assets/misc/cc.png"
    "G61" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A62" = "Ardea pacifica"
    "B62" = "https://www.xeno-canto.org/439287/download"
    "C62" = "Chong Swamp, Queensland, Australia"
    "D62" = "Marc Anderson"
    "E62" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F62" = "assets/misc/cc.png"
    "G62" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A63" = "Ardenna carneipes"
    "B63" = "https://www.xeno-canto.org/424726/download"
    "C63" = "Lady Alice Island, Whangarei District, New Zealand"
    "D63" = "David Boyle"
    "E63" = "https://www.xeno-canto.org/contributor/YHOCFQHBDL"
    "F63" = "assets/misc/cc.png"
    "G63" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A64" = "Ardenna pacifica"
    "B64" = "https://www.xeno-canto.org/113088/download"
    "C64" = "Wilson Island, Queensland, Australia"
    "D64" = "Simon Elliott"
    "E64" = "https://www.xeno-canto.org/contributor/FFFADKCCII"
    "F64" = "assets/misc/cc.png"
    "G64" = "https://creativecommons.org/licenses/by-nc-sa/3.0/"
    "A65" = "Ardenna tenuirostris"
    "B65" = "https://www.xeno-canto.org/177452/download"
    "C65" = "Phillip Island, Victoria, Australia"
    "D65" = "Nick Talbot"
    "E65" = "https://www.xeno-canto.org/contributor/CCUCXWCPSW"
    "F65" = "assets/misc/cc.png"
    "G65" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A66" = "Arenaria interpres"
    "B66" = "https://www.xeno-canto.org/435117/download"
    "C66" = "Snaefellsness, Iceland"
    "D66" = "Patrick Franke"
    "E66" = "https://www.xeno-canto.org/435117/download"
    "F66" = "assets/misc/cc.png"
    "G66" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A67" = "Arenaria interpres"
    "B67" = "https://www.xeno-canto.org/311155/download"
    "C67" = "Bundala Wetland, Sri Lanka"
    "D67" = "Peter Boesman"
    "E67" = "https://www.xeno-canto.org/contributor/OOECIWCSWV"
    "F67" = "assets/misc/cc.png"
    "G67" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A68" = "Artamus cinereus"
    "B68" = "https://www.xeno-canto.org/428464/download"
    "C68" = "Roebuck Plains Station, Western Australia, Australia"
    "D68" = "Nigel Jackett"
    "E68" = "https://www.xeno-canto.org/contributor/KXKBPMRFTY"
    "F68" = "assets/misc/cc.png"
    "G68" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A69" = "Artamus cinereus"
    "B69" = "https://www.xeno-canto.org/334914/download"
    "C69" = "Tibooburra, New South Wales, Australia"
    "D69" = "Marc Anderson"
    "E69" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F69" = "assets/misc/cc.png"
    "G69" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A70" = "Artamus cyanopterus"
    "B70" = "https://www.xeno-canto.org/354593/download"
    "C70" = "Bowra Station, Queensland, Australia"
    "D70" = "Greg McLachlan"
    "E70" = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
    "F70" = "assets/misc/cc.png"
    "G70" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A71" = "Artamus leucorynchus"
    "B71" = "https://www.xeno-canto.org/283305/download"
    "C71" = "Tuggerah, New South Wales, Australia"
    "D71" = "Marc Anderson"
    "E71" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F71" = "assets/misc/cc.png"
    "G71" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A72" = "Artamus leucorynchus"
    "B72" = "https://www.xeno-canto.org/283304/download"
    "C72" = "Tuggerah, New South Wales, Australia"
    "D72" = "Marc Anderson"
    "E72" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F72" = "assets/misc/cc.png"
    "G72" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A73" = "Artamus minor"
    "B73" = "https://www.xeno-canto.org/107836/download"
    "C73" = "Kalbarri National Park, Western Australia, Australia"
    "D73" = "Matthias Feuersenger"
    "E73" = "https://www.xeno-canto.org/contributor/HBPYQXTJEV"
    "F73" = "assets/misc/cc.png"
    "G73" = "https://creativecommons.org/licenses/by-nc-nd/2.5/"
    "A74" = "Artamus personatus"
    "B74" = "https://www.xeno-canto.org/287053/download"
    "C74" = "Bowra Station, Queensland, Australia"
    "D74" = "Marc Anderson"
    "E74" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F74" = "assets/misc/cc.png"
    "G74" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A75" = "Artamus superciliosus"
    "B75" = "https://www.xeno-canto.org/407726/download"
    "C75" = "Mt Ida, Victoria, Australia"
    "D75" = "Frank Lambert"
    "E75" = "https://www.xeno-canto.org/contributor/YTUXOCTUEM"
    "F75" = "assets/misc/cc.png"
    "G75" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A76" = "Artamus superciliosus"
    "B76" = "https://www.xeno-canto.org/393405/download"
    "C76" = "Eubalong, New South Wales, Australia"
    "D76" = "Greg McLachlan"
    "E76" = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
    "F76" = "assets/misc/cc.png"
    "G76" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
    "A77" = "Artamus superciliosus"
    "B77" = "https://www.xeno-canto.org/389388/download"
    "C77" = "Pitt Town Lagoon, New South Wales, Australia"
    "D77" = "Marc Anderson"
    "E77" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F77" = "assets/misc/cc.png"
    "G77" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A78" = "Aviceda subcristata"
    "B78" = "https://www.xeno-canto.org/438827/download"
    "C78" = "Oyala-Thumotang National Park, Queensland, Australia"
    "D78" = "Marc Anderson"
    "E78" = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
    "F78" = "assets/misc/cc.png"
    "G78" = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
    "A79" = "Aviceda subcristata"
    "B79" = "https://www.xeno-canto.org/409476/download"
    "C79" = "Kabupaten Raja Ampat, Papua Barat, Indonesia"
    "D79" = "Ross Gallardy"
    "E79" = "https://www.xeno-canto.org/contributor/FNIOJOZADD"
    "F79" = "assets/misc/cc.png"
    "G79" = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
